$d = $word.ActiveDocument

# The paragraph contains pseudo-markup text that reads (visibly):
#   ... se regorge &amp; <del>en</del><lb/> ...
# In the OOXML the "en" inside <del>...</del> is split across two
# adjacent runs ("e" and "n", same rPr: strike=0 / color=000000 / rtl=0).
# We need to merge them into a single run whose text is "en".
#
# Locate the unique anchor immediately before the "en" text so we can
# build a tight 2-character Range over exactly those two runs, without
# touching the (differently formatted) "<del>" / "</del>" markup runs
# around it.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "se regorge &amp; <del>", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the anchor text before the target run pair."
}

$anchor.Collapse(0)
$start = $anchor.Start
$target = $d.Range($start, $start + 2)

if ($target.Text -ne "en") {
    throw "Unexpected target text: [$($target.Text)]"
}

# A plain "set the same text" is a no-op for this engine (no visible
# text change => no run restructuring), so force a genuine text change
# first, then restore the final text; that collapses the run pair into
# a single run carrying the original (shared) run formatting.
$target.Text = "__"
$target = $d.Range($start, $start + 2)
$target.Text = "en"
